# JS ES6: Use Destructuring Assignment to Assign Variables from Objects
# Marks the next 10 freeCodeCamp ES6 challenges (rows 113-122) as "Passed"
# (new text + completion date) and highlights the following block of
# still-outstanding challenges (rows 123-132) in bold so it's clear which
# lesson is now "up next".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared strings / completion date (date serial 43436 = 2018-12-02)
$doneDate = 43436

$passedText = @{
    113 = "PassedDeclare a Read-Only Variable with the const Keyword"
    114 = "PassedMutate an Array Declared with const"
    115 = "PassedPrevent Object Mutation"
    116 = "PassedUse Arrow Functions to Write Concise Anonymous Functions"
    117 = "PassedWrite Arrow Functions with Parameters"
    118 = "PassedWrite Higher Order Arrow Functions"
    119 = "PassedSet Default Parameters for Your Functions"
    120 = "PassedUse the Rest Operator with Function Parameters"
    121 = "PassedUse the Spread Operator to Evaluate Arrays In-Place"
    122 = "PassedUse Destructuring Assignment to Assign Variables from Objects"
}

foreach ($row in 113..122) {
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value2 = $passedText[$row]
    # Drop the "not passed" highlight fill now that the challenge is done.
    $dCell.ClearFormats()
    $dCell.Interior.Pattern = -4142

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value2 = $doneDate
}

# The next block of challenges (still outstanding) gets bolded so it stands
# out as the new "current" section.
foreach ($row in 123..132) {
    $ws.Cells.Item($row, 4).Font.Bold = $true
}

# Move the frozen-pane view / selection to the newly-updated rows.
$ws.Activate()
$sel = $ws.Range("D121")
$sel.Select()
